$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Add new row 51 with the new mail log entry
$ws.Range("A51").Value = "Kan mijn wachtwoord niet resetten"
$ws.Range("B51").Value = "mailmind.test@zohomail.eu"
$ws.Range("C51").Value = "Ik krijg geen e-mail bij wachtwoord resetten."
$ws.Range("D51").Value = "IT / Technisch probleem"
$ws.Range("E51").Value = "Beste klant,`nBedankt voor het melden van dit probleem. Om u verder te kunnen helpen, hebben we enkele aanvullende gegevens nodig. Kunt u ons uw gebruikersnaam of het e-mailadres waarmee u probeerde uw wachtwoord opnieuw in te stellen, geven? Op die manier kunnen we het probleem gericht onderzoeken en oplossen.`nAlvast bedankt voor uw medewerking.`nMet vriendelijke groet,`n[Naam bedrijf] - E-mailassistent"
$ws.Range("F51").Value = "2025-06-22 21:58:13"
$ws.Range("G51").Value = "Ja"

# Reset row height to the sheet default (remove autosize artifact from the
# multi-line text so the row matches the other plain data rows)
$ws.Rows.Item(51).AutoFit()

# Extend the conditional formatting ranges to include the new row
$fcD = $ws.Range("D2:D50").FormatConditions
$fcD.Item(1).ModifyAppliesToRange($ws.Range("D2:D51"))

$fcG = $ws.Range("G2:G50").FormatConditions
$fcG.Item(1).ModifyAppliesToRange($ws.Range("G2:G51"))

# Update the Dashboard summary count for "IT / Technisch probleem"
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 9
